# Apply the "Updated cryptos list" data refresh to Sheet1.
# Columns: A=index(unchanged) B=Coin C=Link D=Price E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.972.43"
$ws.Range("E2").Value = "  +0.41%  "

# Row 3
$ws.Range("D3").Value = "2.548.55"
$ws.Range("E3").Value = "  +0.54%  "

# Row 4
$ws.Range("E4").Value = "  -0.19%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.66"
$ws.Range("E5").Value = "  +1.85%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.99"
$ws.Range("E6").Value = "  +6.57%  "

# Row 7
$ws.Range("E7").Value = "  +1.03%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.549"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.56"
$ws.Range("E10").Value = "  +2.14%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0836"
$ws.Range("E11").Value = "  +4.32%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.115"
$ws.Range("E12").Value = "  +2.00%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.66"
$ws.Range("E13").Value = "  -0.42%  "

# Row 14
$ws.Range("D14").Value = "2.941.11"
$ws.Range("E14").Value = "  +0.58%  "

# Row 15
$ws.Range("D15").Value = "2.599.09"
$ws.Range("E15").Value = "  +2.50%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.13"
$ws.Range("E16").Value = "  +7.49%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.876"
$ws.Range("E17").Value = "  +0.99%  "

# Row 18
$ws.Range("D18").Value = "42.989.23"
$ws.Range("E18").Value = "  +0.31%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.90"
$ws.Range("E19").Value = "  +6.84%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0999"
$ws.Range("E20").Value = "  +1.30%  "

# Row 21
$ws.Range("E21").Value = "  +0.82%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.94"
$ws.Range("E22").Value = "  +0.67%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "255.60"
$ws.Range("E23").Value = "  -0.32%  "

# Row 24
$ws.Range("E24").Value = "  +1.77%  "

# Row 25
$ws.Range("E25").Value = "  -0.71%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.05"
$ws.Range("E26").Value = "  -4.18%  "

# Row 27
$ws.Range("E27").Value = "  -0.09%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.25"
$ws.Range("E28").Value = "  +2.53%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.75"
$ws.Range("E29").Value = "  +2.03%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.21"
$ws.Range("E30").Value = "  +4.91%  "

# Row 31
$ws.Range("E31").Value = "  -1.43%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.76"
$ws.Range("E32").Value = "  +4.26%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.87"
$ws.Range("E33").Value = "  +18.40%  "

# Row 34
$ws.Range("E34").Value = "  -0.93%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0804"
$ws.Range("E35").Value = "  +1.50%  "

# Row 36
$ws.Range("E36").Value = "  -2.05%  "

# Row 38
$ws.Range("E38").Value = "  +2.62%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.35"
$ws.Range("E39").Value = "  +5.89%  "

# Row 40
$ws.Range("E40").Value = "  +0.46%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.92"
$ws.Range("E41").Value = "  +0.71%  "

# Row 42
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.44"
$ws.Range("E42").Value = "  +1.62%  "

# Row 43
$ws.Range("E43").Value = "  +28.33%  "

# Row 44
$ws.Range("D44").Value = "2.105.26"
$ws.Range("E44").Value = "  +0.93%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0306"
$ws.Range("E45").Value = "  -1.25%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.99"
$ws.Range("E47").Value = "  +3.56%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.99"
$ws.Range("E48").Value = "  +0.66%  "

# Row 49
$ws.Range("D49").Value = "2.799.29"
$ws.Range("E49").Value = "  +0.58%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.78"
$ws.Range("E50").Value = "  +9.10%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "103.58"
$ws.Range("E51").Value = "  -0.34%  "
